# Apply the edit described by the commit:
# "Criação Notebook com Spark" - clean up the AnaliseExploratoria sheet:
#   - clear the now-unused TIPO (String/INT) values in column C, rows 3-12
#   - remove the whole row for the 'IN_MATERIAL_ESPECÍFICO' variable (row 72)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AnaliseExploratoria")
$ws.Activate()

# Clear the TIPO column contents for rows 3 to 12 (keep formatting/styles).
$ws.Range("C3:C12").ClearContents()

# Delete the entire row that documents 'IN_MATERIAL_ESPECÍFICO' (row 72),
# shifting all subsequent rows up by one.
$ws.Rows.Item(72).Delete()

# Reflect the resulting active cell selection.
$ws.Range("B13").Select() | Out-Null
